$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adjust column A width (target stored width 13.7109375 character-units;
# the host's column-width model quantizes to 1/6-character steps, so feed
# it the ColumnWidth value whose rounded pixel width lands closest to the
# target: 12.8 -> stored width 13.666666... is the nearest reachable value)
$ws.Columns.Item(1).ColumnWidth = 12.8

# Update values for rows 1-3 (rows 4-5 unchanged)
$ws.Range("A1").Value = 0.030114300841754746
$ws.Range("B1").Value = -0.030114300853086629

$ws.Range("A2").Value = 0.01339667561690793
$ws.Range("B2").Value = -0.013396675629906367

$ws.Range("A3").Value = -0.039365022973392207
$ws.Range("B3").Value = 0.039365022935584061
